# Phase 4.4: Fix remaining validation errors
# - Flip "Experimental" flag from true to false
# - Update the "Date" metadata value to the new timestamp
#
# Both values live on the "Metadata" worksheet, which holds a simple
# Property/Value table: column A has the property name, column B has the
# value. Row 7 = Experimental, Row 8 = Date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental: true -> false
# Leading apostrophe forces Excel to keep this as literal text "false"
# instead of auto-coercing it to the Boolean FALSE (same as typing it
# into the Excel UI).
$ws.Range("B7").Value = "'false"

# Date: refreshed publication timestamp
$ws.Range("B8").Value = "2025-10-03T16:37:46+01:00"
